$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell text updates (coinranking.com crypto price/volume snapshot refresh).
# For cells whose new text looks like a plain number (e.g. "363.30"), force
# text formatting first so Excel does not silently normalize it to 363.3.

$ws.Range("D2").Value = "67.797.02"
$ws.Range("E2").Value = "  +0.73%  "
$ws.Range("D3").Value = "2.499.71"
$ws.Range("E3").Value = "  -1.86%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "592.58"
$ws.Range("E5").Value = "  +0.35%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.39"
$ws.Range("E6").Value = "  -0.17%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.527"
$ws.Range("E8").Value = "  -0.66%  "
$ws.Range("D9").Value = "2.499.65"
$ws.Range("E9").Value = "  -1.83%  "
$ws.Range("E10").Value = "  -0.27%  "
$ws.Range("E11").Value = "  +2.14%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.08"
$ws.Range("E12").Value = "  -1.69%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.342"
$ws.Range("E13").Value = "  -2.75%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.29"
$ws.Range("E14").Value = "  -2.87%  "
$ws.Range("D15").Value = "2.955.60"
$ws.Range("E15").Value = "  -1.63%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000177"
$ws.Range("E16").Value = "  -1.10%  "
$ws.Range("D17").Value = "67.709.40"
$ws.Range("E17").Value = "  +0.76%  "
$ws.Range("D18").Value = "2.480.68"
$ws.Range("E18").Value = "  -2.59%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.79"
$ws.Range("E19").Value = "  +3.12%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.96"
$ws.Range("E20").Value = "  -1.12%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "363.30"
$ws.Range("E21").Value = "  +1.86%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.13"
$ws.Range("E22").Value = "  -2.26%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.56"
$ws.Range("E23").Value = "  -2.63%  "
$ws.Range("E24").Value = "  +0.00%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "71.19"
$ws.Range("E25").Value = "  +1.55%  "
$ws.Range("E26").Value = "  -5.65%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.83"
$ws.Range("E27").Value = "  -2.27%  "
$ws.Range("E28").Value = "  -0.21%  "
$ws.Range("D29").Value = "2.620.80"
$ws.Range("E29").Value = "  -1.78%  "
$ws.Range("D30").Value = "0.0₃0968"
$ws.Range("E30").Value = "  -3.10%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "532.20"
$ws.Range("E31").Value = "  -0.74%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.24"
$ws.Range("E32").Value = "  -0.01%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.88"
$ws.Range("E33").Value = "  +0.79%  "
$ws.Range("E34").Value = "  -3.57%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.127"
$ws.Range("E36").Value = "  -3.94%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "158.66"
$ws.Range("E37").Value = "  +0.98%  "
$ws.Range("E38").Value = "  -3.46%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.61"
$ws.Range("E39").Value = "  -1.15%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.65"
$ws.Range("E40").Value = "  +1.06%  "
$ws.Range("E41").Value = "  -1.56%  "
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.12"
$ws.Range("E42").Value = "  -1.77%  "
$ws.Range("B43").Value = "PolygonEcosystemToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.348"
$ws.Range("E43").Value = "  -2.55%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.999"
$ws.Range("E44").Value = "  -0.14%  "
$ws.Range("E45").Value = "  -2.53%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "145.51"
$ws.Range("E46").Value = "  -3.72%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.69"
$ws.Range("E47").Value = "  -0.98%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.549"
$ws.Range("E48").Value = "  -2.86%  "
$ws.Range("E49").Value = "  -3.69%  "
$ws.Range("E50").Value = "  -1.32%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0750"
$ws.Range("E51").Value = "  -1.47%  "
